# MaxFatLoss.xlsx - "Working on edge conditions of low body fat and losses"
#
# 1) B2 (%Fat) drops from 8.2 to 7.6 -- this ripples through B3/F3 (formulas,
#    recalculated automatically).
# 2) E14's formula gains a low-body-fat fallback branch: instead of clamping
#    to F3 (Loss Fat / unlimited) when F3 <= E13, it now falls back to
#    E8-F3 (Fat available minus Loss Fat) -- which ripples into B14/E15/F15/G15.
# 3) The active selection moves from G15 to E14 (the cell being edited).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New %Fat input value.
$ws.Range("B2").Value = 7.6

# New edge-case formula for "Loss Fat (selected)".
$ws.Range("E14").Formula = "=IF(F3>E13,E13,E8-F3)"

# Selection follows the edit onto E14 (was G15).
$ws.Range("E14").Select()
